$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Trf"
$ws.Range("C2").Value = "Tfrc"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 12.43612633333333
$ws.Range("H2").Value = 37.308379
$ws.Range("I2").Value = 0.2323717069953836
$ws.Range("J2").Value = 0.2323717069953836
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 27.631265
$ws.Range("N2").Value = 82.893795
$ws.Range("O2").Value = 0.5169298047761691
$ws.Range("P2").Value = 0.516929804776169
$ws.Range("Q2").Value = 343.6259022898117
$ws.Range("R2").Value = 3092.633120608305
$ws.Range("S2").Value = 0.1201198611326288
$ws.Range("T2").Value = 0.1201198611326288

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Trf"
$ws.Range("C3").Value = "Tfrc"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 12.43612633333333
$ws.Range("H3").Value = 37.308379
$ws.Range("I3").Value = 0.2323717069953836
$ws.Range("J3").Value = 0.2323717069953836
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.171031666666667
$ws.Range("N3").Value = 6.513095000000001
$ws.Range("O3").Value = 0.0406159824971054
$ws.Range("P3").Value = 0.0406159824971054
$ws.Range("Q3").Value = 26.99922408033389
$ws.Range("R3").Value = 242.993016723005
$ws.Range("S3").Value = 0.009438005184147006
$ws.Range("T3").Value = 0.009438005184147006

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Trf"
$ws.Range("C4").Value = "Tfrc"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 12.43612633333333
$ws.Range("H4").Value = 37.308379
$ws.Range("I4").Value = 0.2323717069953836
$ws.Range("J4").Value = 0.2323717069953836
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 23.65034766666666
$ws.Range("N4").Value = 70.951043
$ws.Range("O4").Value = 0.4424542127267255
$ws.Range("P4").Value = 0.4424542127267255
$ws.Range("Q4").Value = 294.1187114099219
$ws.Range("R4").Value = 2647.068402689297
$ws.Range("S4").Value = 0.1028138406786078
$ws.Range("T4").Value = 0.1028138406786078

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Trf"
$ws.Range("C5").Value = "Tfrc"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 39.70924033333333
$ws.Range("H5").Value = 119.127721
$ws.Range("I5").Value = 0.7419757336345224
$ws.Range("J5").Value = 0.7419757336345223
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 27.631265
$ws.Range("N5").Value = 82.893795
$ws.Range("O5").Value = 0.5169298047761691
$ws.Range("P5").Value = 0.516929804776169
$ws.Range("Q5").Value = 1097.216542599022
$ws.Range("R5").Value = 9874.948883391195
$ws.Range("S5").Value = 0.3835493711363485
$ws.Range("T5").Value = 0.3835493711363484

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Trf"
$ws.Range("C6").Value = "Tfrc"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 39.70924033333333
$ws.Range("H6").Value = 119.127721
$ws.Range("I6").Value = 0.7419757336345224
$ws.Range("J6").Value = 0.7419757336345223
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.171031666666667
$ws.Range("N6").Value = 6.513095000000001
$ws.Range("O6").Value = 0.0406159824971054
$ws.Range("P6").Value = 0.0406159824971054
$ws.Range("Q6").Value = 86.2100182229439
$ws.Range("R6").Value = 775.8901640064951
$ws.Range("S6").Value = 0.0301360734105767
$ws.Range("T6").Value = 0.0301360734105767

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Trf"
$ws.Range("C7").Value = "Tfrc"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 39.70924033333333
$ws.Range("H7").Value = 119.127721
$ws.Range("I7").Value = 0.7419757336345224
$ws.Range("J7").Value = 0.7419757336345223
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 23.65034766666666
$ws.Range("N7").Value = 70.951043
$ws.Range("O7").Value = 0.4424542127267255
$ws.Range("P7").Value = 0.4424542127267255
$ws.Range("Q7").Value = 939.1373394625558
$ws.Range("R7").Value = 8452.236055163003
$ws.Range("S7").Value = 0.3282902890875972
$ws.Range("T7").Value = 0.3282902890875972

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Trf"
$ws.Range("C8").Value = "Tfrc"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.37288
$ws.Range("H8").Value = 4.11864
$ws.Range("I8").Value = 0.02565255937009396
$ws.Range("J8").Value = 0.02565255937009396
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 27.631265
$ws.Range("N8").Value = 82.893795
$ws.Range("O8").Value = 0.5169298047761691
$ws.Range("P8").Value = 0.516929804776169
$ws.Range("Q8").Value = 37.9344110932
$ws.Range("R8").Value = 341.4096998388
$ws.Range("S8").Value = 0.01326057250719176
$ws.Range("T8").Value = 0.01326057250719176

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Trf"
$ws.Range("C9").Value = "Tfrc"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.37288
$ws.Range("H9").Value = 4.11864
$ws.Range("I9").Value = 0.02565255937009396
$ws.Range("J9").Value = 0.02565255937009396
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.171031666666667
$ws.Range("N9").Value = 6.513095000000001
$ws.Range("O9").Value = 0.0406159824971054
$ws.Range("P9").Value = 0.0406159824971054
$ws.Range("Q9").Value = 2.980565954533334
$ws.Range("R9").Value = 26.82509359080001
$ws.Range("S9").Value = 0.001041903902381694
$ws.Range("T9").Value = 0.001041903902381694

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Trf"
$ws.Range("C10").Value = "Tfrc"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.37288
$ws.Range("H10").Value = 4.11864
$ws.Range("I10").Value = 0.02565255937009396
$ws.Range("J10").Value = 0.02565255937009396
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 23.65034766666666
$ws.Range("N10").Value = 70.951043
$ws.Range("O10").Value = 0.4424542127267255
$ws.Range("P10").Value = 0.4424542127267255
$ws.Range("Q10").Value = 32.46908930461333
$ws.Range("R10").Value = 292.22180374152
$ws.Range("S10").Value = 0.01135008296052051
$ws.Range("T10").Value = 0.01135008296052051
